$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 109
$ws.Range("C2").Value = "face/face008.jpg"
$ws.Range("D2").Value = "schmecken"
$ws.Range("E2").Value = "face"
$ws.Range("B3").Value = 101
$ws.Range("C3").Value = "face/face014.jpg"
$ws.Range("D3").Value = "runden"
$ws.Range("E3").Value = "face"
$ws.Range("B4").Value = 42
$ws.Range("C4").Value = "house/house008.jpg"
$ws.Range("D4").Value = "pflegen"
$ws.Range("E4").Value = "house"
$ws.Range("B5").Value = 75
$ws.Range("C5").Value = "house/house020.jpg"
$ws.Range("D5").Value = "krachen"
$ws.Range("E5").Value = "house"
$ws.Range("B6").Value = 72
$ws.Range("C6").Value = "house/house001.jpg"
$ws.Range("D6").Value = "bleiben"
$ws.Range("E6").Value = "house"
$ws.Range("B7").Value = 59
$ws.Range("C7").Value = "face/face019.jpg"
$ws.Range("D7").Value = "wenden"
$ws.Range("E7").Value = "face"
$ws.Range("B8").Value = 82
$ws.Range("C8").Value = "house/house030.jpg"
$ws.Range("D8").Value = "laufen"
$ws.Range("E8").Value = "house"
$ws.Range("B9").Value = 19
$ws.Range("C9").Value = "face/face017.jpg"
$ws.Range("D9").Value = "währen"
$ws.Range("E9").Value = "face"
$ws.Range("B10").Value = 105
$ws.Range("C10").Value = "face/face029.jpg"
$ws.Range("D10").Value = "rasen"
$ws.Range("E10").Value = "face"
$ws.Range("B11").Value = 29
$ws.Range("C11").Value = "house/house007.jpg"
$ws.Range("D11").Value = "segeln"
$ws.Range("E11").Value = "house"
$ws.Range("B12").Value = 120
$ws.Range("C12").Value = "face/face007.jpg"
$ws.Range("D12").Value = "mieten"
$ws.Range("E12").Value = "face"
$ws.Range("B13").Value = 23
$ws.Range("C13").Value = "house/house022.jpg"
$ws.Range("D13").Value = "loben"
$ws.Range("E13").Value = "house"
$ws.Range("B14").Value = 37
$ws.Range("C14").Value = "house/house015.jpg"
$ws.Range("D14").Value = "jubeln"
$ws.Range("E14").Value = "house"
$ws.Range("B15").Value = 74
$ws.Range("C15").Value = "house/house009.jpg"
$ws.Range("D15").Value = "antun"
$ws.Range("E15").Value = "house"
$ws.Range("B16").Value = 60
$ws.Range("C16").Value = "face/face012.jpg"
$ws.Range("D16").Value = "füttern"
$ws.Range("E16").Value = "face"
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = "house/house010.jpg"
$ws.Range("D17").Value = "wiegen"
$ws.Range("E17").Value = "house"
$ws.Range("B18").Value = 117
$ws.Range("C18").Value = "face/face003.jpg"
$ws.Range("D18").Value = "regnen"
$ws.Range("E18").Value = "face"
$ws.Range("B19").Value = 56
$ws.Range("C19").Value = "face/face022.jpg"
$ws.Range("D19").Value = "spielen"
$ws.Range("E19").Value = "face"
$ws.Range("B20").Value = 125
$ws.Range("C20").Value = "face/face025.jpg"
$ws.Range("D20").Value = "kaufen"
$ws.Range("E20").Value = "face"
$ws.Range("B21").Value = 30
$ws.Range("C21").Value = "face/face030.jpg"
$ws.Range("D21").Value = "scheitern"
$ws.Range("E21").Value = "face"
$ws.Range("B22").Value = 87
$ws.Range("C22").Value = "house/house006.jpg"
$ws.Range("D22").Value = "gelten"
$ws.Range("E22").Value = "house"
$ws.Range("B23").Value = 121
$ws.Range("C23").Value = "house/house026.jpg"
$ws.Range("D23").Value = "fliegen"
$ws.Range("E23").Value = "house"
$ws.Range("B24").Value = 112
$ws.Range("C24").Value = "face/face009.jpg"
$ws.Range("D24").Value = "stechen"
$ws.Range("E24").Value = "face"
$ws.Range("B25").Value = 80
$ws.Range("C25").Value = "house/house021.jpg"
$ws.Range("D25").Value = "formen"
$ws.Range("E25").Value = "house"
$ws.Range("B26").Value = 51
$ws.Range("C26").Value = "house/house028.jpg"
$ws.Range("D26").Value = "schenken"
$ws.Range("E26").Value = "house"
$ws.Range("B27").Value = 90
$ws.Range("C27").Value = "face/face018.jpg"
$ws.Range("D27").Value = "fühlen"
$ws.Range("E27").Value = "face"
$ws.Range("B28").Value = 98
$ws.Range("C28").Value = "face/face006.jpg"
$ws.Range("D28").Value = "langen"
$ws.Range("E28").Value = "face"
$ws.Range("B29").Value = 48
$ws.Range("C29").Value = "house/house005.jpg"
$ws.Range("D29").Value = "starten"
$ws.Range("E29").Value = "house"
$ws.Range("B30").Value = 97
$ws.Range("C30").Value = "face/face027.jpg"
$ws.Range("D30").Value = "rücken"
$ws.Range("E30").Value = "face"
$ws.Range("B31").Value = 76
$ws.Range("C31").Value = "house/house016.jpg"
$ws.Range("D31").Value = "klappen"
$ws.Range("E31").Value = "house"
$ws.Range("B32").Value = 63
$ws.Range("C32").Value = "house/house003.jpg"
$ws.Range("D32").Value = "husten"
$ws.Range("E32").Value = "house"
$ws.Range("B33").Value = 99
$ws.Range("C33").Value = "face/face004.jpg"
$ws.Range("D33").Value = "nehmen"
$ws.Range("E33").Value = "face"
